# Append two new daily rows (2025-12-06 / serial 45997) to Sheet1,
# one for each station, mirroring the existing row layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 - 四方坪站
$ws.Cells.Item(12, 1).Value = 45997
$ws.Cells.Item(12, 2).Value = "四方坪站"
$ws.Cells.Item(12, 3).Value = 10156.379999999999
$ws.Cells.Item(12, 4).Value = 8504.85
$ws.Cells.Item(12, 5).Value = 3371.76
$ws.Cells.Item(12, 6).Value = 436

# Row 13 - 高岭站
$ws.Cells.Item(13, 1).Value = 45997
$ws.Cells.Item(13, 2).Value = "高岭站"
$ws.Cells.Item(13, 3).Value = 4756.08
$ws.Cells.Item(13, 4).Value = 4044.68
$ws.Cells.Item(13, 5).Value = 1143.9100000000001
$ws.Cells.Item(13, 6).Value = 170

# Match the author's final selection/view state as closely as possible
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I16").Select()
